# "Export all in words implemented"
#
# Adds an "Exams" results table (ID / Name / Time / Status / Score /
# Examiner header row) right after the existing "Exams" heading
# paragraph, before the section properties at the end of the body.

$d = $word.ActiveDocument

# Insert the table at the very end of the document's content, i.e.
# immediately after the last paragraph ("Exams") and before the
# sectPr - this keeps the "Exams" paragraph intact instead of
# splitting it.
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$tbl = $d.Tables.Add($insertionPoint, 1, 6)

# Match the tblLook flags used by the target table: heading row +
# first column emphasised, row banding on, column banding off.
$tbl.ApplyStyleHeadingRows = $true
$tbl.ApplyStyleFirstColumn = $true
$tbl.ApplyStyleLastColumn = $false
$tbl.ApplyStyleLastRow = $false
$tbl.ApplyStyleRowBands = $true
$tbl.ApplyStyleColumnBands = $false

# Every column is 1440 twips (1 inch) wide == 72 points.
for ($i = 1; $i -le $tbl.Columns.Count; $i++) {
    $tbl.Columns($i).Width = 72
}

# Fill in the header row, resetting each cell's paragraph style back
# to Normal so it doesn't inherit the "Exams" heading's style.
$headers = @("ID", "Name", "Time", "Status", "Score", "Examiner")
for ($i = 1; $i -le 6; $i++) {
    $cell = $tbl.Cell(1, $i)
    $cell.Range.Style = "Normal"
    $cell.Range.Text = $headers[$i - 1]
}
